$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update sheet (tab) name to reflect new "through" date
$ws.Name = "Through 2022-08-21"

# Update label for August row
$ws.Range("A9").Value = "August (through 08-21)"

# Update August (row 9) figures
$ws.Range("C9").Value = 49
$ws.Range("D9").Value = 55
$ws.Range("E9").Value = 32
$ws.Range("G9").Value = 128
$ws.Range("H9").Value = 110
$ws.Range("I9").Value = 123

# Update Total (row 10) figures
$ws.Range("C10").Value = 351
$ws.Range("D10").Value = 520
$ws.Range("E10").Value = 457
$ws.Range("G10").Value = 749
$ws.Range("H10").Value = 1020
$ws.Range("I10").Value = 1094
